$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Labels for rows 2..45 (column A), values move to column B
$labels = @(
    "rendimento médio real habitual do trabalho principal(r$/mês)",
    "rendimento médio real habitual de todos os trabalhos(r$/mês)",
    "índice gini",
    "rendimento-hora médio real habitual de todos os trabalhos(r$/hora)",
    "rendimento-hora médio real habitual do trabalho principal(r$/hora)",
    "taxa de formalização",
    "sem instrução ou fundamental incompleto",
    "taxa composta de subutilização",
    "60 a 69 anos",
    "30 a 44 anos",
    "15 a 29 anos",
    "ensino fundamental completo ou médio incompleto",
    "domicílio próprio - pagando",
    "taxa de desocupação",
    "domicílio cedido por empregador",
    "taxa de analfabetismo",
    "população ocupada",
    "população desocupada",
    "domicílio próprio - já pago",
    "70 anos ou mais",
    "domicílio alugado",
    "outra forma",
    "população na força de trabalho",
    "taxa total mortalidade",
    "população ocupada em trabalhos formais",
    "número de beneficiários de plano de saúde",
    "total pessoas ocupadas(1 000 pessoas)",
    "ensino médio completo ou superior incompleto",
    "população em idade de trabalhar",
    "45 a 59 anos",
    "população na força de trabalho potencial",
    "população subutilizada",
    "domicílio cedido por familiar",
    "domicílio cedido de outra forma",
    "taxa de participação",
    "saneamento basico total(1 000 pessoas)",
    "total pessoas por condição de ocupação a domicílio(1 000 pessoas)",
    "0 a 14 anos",
    "número mensal médio de leitos de internação (total)",
    "ensino superior completo",
    "população",
    "nível de ocupação",
    "proporção de pessoas com acesso simultâneo aos três serviços de saneamento básico(%)",
    "60 anos ou mais"
)

$values = @(
    0.2408144716046481,
    0.2138535408764226,
    0.1707172627269967,
    0.1107999724384888,
    0.1001127194766462,
    0.09436755344224074,
    -0.07104905201899941,
    -0.06795054771273268,
    0.06332571302153378,
    -0.05712421373858256,
    -0.05706900969261897,
    -0.05510071196494273,
    0.04275217488633623,
    -0.04234924181250599,
    -0.04059600860080881,
    -0.03970962527386596,
    -0.03769138807096256,
    0.03623613648447817,
    -0.03387290178837965,
    -0.03380587008850378,
    0.03290344154141278,
    0.03222123472848953,
    -0.02982483317104387,
    -0.02870304499982512,
    -0.02855980567584415,
    -0.02762849389152877,
    0.02391757285780738,
    -0.02303298864907466,
    -0.02208524538301557,
    0.02021892966958827,
    0.0164726686880224,
    0.0161930030681766,
    0.01553472395333394,
    -0.01553355023422243,
    -0.01323612318788103,
    -0.01113816392471648,
    -0.01113816392471648,
    0.01029084974979556,
    -0.009798692804159598,
    0.006998962984198775,
    -0.005287046429947873,
    -0.003604956937174886,
    -0.003443971416822456,
    -0.0014919028399486
)

# Header row: A1 becomes "variavel", B1 becomes "impacto" (copy A1's header style to B1)
$ws.Range("A1").Value = "variavel"
$ws.Range("B1").Value = "impacto"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: labels in column A, numeric impact values moved to column B
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
